$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3699
$ws1.Range("F5").Value = 3699
$ws1.Range("F6").Value = 280
$ws1.Range("F7").Value = 5229
$ws1.Range("F9").Value = 396
$ws1.Range("F11").Value = 723
$ws1.Range("F14").Value = 41
$ws1.Range("F16").Value = 338
$ws1.Range("F19").Value = 165
$ws1.Range("F22").Value = 5987
$ws1.Range("F26").Value = 6298
$ws1.Range("F34").Value = 131
$ws1.Range("F41").Value = 1079
$ws1.Range("F42").Value = 2046

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3699
$ws4.Range("F8").Value = 3699
$ws4.Range("F9").Value = 280
$ws4.Range("F10").Value = 5229
$ws4.Range("F12").Value = 396
$ws4.Range("F14").Value = 723
$ws4.Range("F17").Value = 41
$ws4.Range("F19").Value = 338
$ws4.Range("F23").Value = 165
$ws4.Range("F26").Value = 5987
$ws4.Range("F30").Value = 6298
$ws4.Range("F39").Value = 131
$ws4.Range("F46").Value = 1079
$ws4.Range("F48").Value = 2046

$wb.Save()
